$d = $word.ActiveDocument

$replacements = @(
    @("966×3=2898", "436×2=872"),
    @("648×9=5832", "937×3=2811"),
    @("965×5=4825", "303×4=1212"),
    @("932×8=7456", "736×9=6624"),
    @("672×6=4032", "638×2=1276"),
    @("344×3=1032", "603×6=3618"),
    @("381×9=3429", "741×2=1482"),
    @("209×8=1672", "486×2=972"),
    @("142×9=1278", "587×4=2348"),
    @("884×6=5304", "221×5=1105"),
    @("474×8=3792", "875×8=7000"),
    @("514×4=2056", "805×5=4025"),
    @("273×6=1638", "925×3=2775"),
    @("665×8=5320", "719×5=3595"),
    @("880×9=7920", "802×4=3208"),
    @("706×9=6354", "843×7=5901"),
    @("717×3=2151", "801×2=1602"),
    @("528×3=1584", "109×7=763"),
    @("119×2=238", "798×3=2394"),
    @("768×6=4608", "708×7=4956"),
    @("236×9=2124", "743×2=1486"),
    @("127×6=762", "343×2=686"),
    @("785×6=4710", "683×6=4098"),
    @("650×8=5200", "547×8=4376"),
    @("768×5=3840", "384×9=3456"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements."